$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '66.816.24'
$ws.Range("E2").Value = '  -3.92%  '
$ws.Range("D3").Value = '3.340.62'
$ws.Range("E3").Value = '  -1.12%  '
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '573.52'
$ws.Range("E5").Value = '  -3.27%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '181.42'
$ws.Range("E6").Value = '  -5.30%  '
$ws.Range("E7").Value = '  -0.02%  '
$ws.Range("E8").Value = '  -1.62%  '
$ws.Range("E9").Value = '  -3.29%  '
$ws.Range("E10").Value = '  -1.58%  '
$ws.Range("E11").Value = '  -4.24%  '
$ws.Range("D12").Value = '3.920.71'
$ws.Range("E12").Value = '  -0.99%  '
$ws.Range("E13").Value = '  -1.75%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '27.12'
$ws.Range("E14").Value = '  -5.29%  '
$ws.Range("D15").Value = '66.879.03'
$ws.Range("E15").Value = '  -3.85%  '
$ws.Range("E16").Value = '  -2.52%  '
$ws.Range("D17").Value = '3.334.56'
$ws.Range("E17").Value = '  -1.05%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '437.53'
$ws.Range("E18").Value = '  -2.42%  '
$ws.Range("B19").Value = 'Polkadot'
$ws.Range("C19").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '5.69'
$ws.Range("E19").Value = '  -2.78%  '
$ws.Range("B20").Value = 'Chainlink'
$ws.Range("C20").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '13.61'
$ws.Range("E20").Value = '  -1.66%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '7.61'
$ws.Range("E21").Value = '  -2.55%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '73.88'
$ws.Range("E22").Value = '  -1.10%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.00'
$ws.Range("E23").Value = '  -0.15%  '
$ws.Range("E24").Value = '  -0.20%  '
$ws.Range("E25").Value = '  -3.83%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.190'
$ws.Range("E26").Value = '  -0.25%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.04'
$ws.Range("E27").Value = '  -5.06%  '
$ws.Range("E28").Value = '  +0.11%  '
$ws.Range("E29").Value = '  -1.49%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '22.84'
$ws.Range("E30").Value = '  -2.16%  '
$ws.Range("E31").Value = '  -6.07%  '
$ws.Range("E32").Value = '  +0.01%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.80'
$ws.Range("E33").Value = '  -2.98%  '
$ws.Range("E34").Value = '  -3.82%  '
$ws.Range("E35").Value = '  -3.08%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '161.95'
$ws.Range("E36").Value = '  -2.13%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '27.91'
$ws.Range("E37").Value = '  +2.17%  '
$ws.Range("E38").Value = '  -5.15%  '
$ws.Range("D39").Value = '2.832.17'
$ws.Range("E39").Value = '  +3.22%  '
$ws.Range("E40").Value = '  -3.33%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '4.45'
$ws.Range("E41").Value = '  -3.43%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '6.30'
$ws.Range("E42").Value = '  -3.20%  '
$ws.Range("B43").Value = 'OKB'
$ws.Range("C43").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '40.12'
$ws.Range("E43").Value = '  -1.60%  '
$ws.Range("B44").Value = 'Hedera'
$ws.Range("C44").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0671'
$ws.Range("E44").Value = '  -2.89%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '24.57'
$ws.Range("E45").Value = '  -4.44%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.36'
$ws.Range("E46").Value = '  -7.03%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '324.24'
$ws.Range("E47").Value = '  -5.45%  '
$ws.Range("E48").Value = '  -4.20%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.986'
$ws.Range("E49").Value = '  -3.55%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '31.05'
$ws.Range("E50").Value = '  -5.85%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '6.17'
$ws.Range("E51").Value = '  -2.90%  '
